# Add a new "Height" property row to the NPC "Property" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$newRow = 44

$ws.Cells.Item($newRow, 1).Value  = "Height"     # Id
$ws.Cells.Item($newRow, 2).Value  = "float"       # Type
$ws.Cells.Item($newRow, 3).Value  = $false        # Public
$ws.Cells.Item($newRow, 4).Value  = $false        # Private
$ws.Cells.Item($newRow, 5).Value  = $false        # Save
$ws.Cells.Item($newRow, 6).Value  = $true         # View
$ws.Cells.Item($newRow, 7).Value  = 0             # Index
$ws.Cells.Item($newRow, 8).Value  = 0             # SaveInterval
$ws.Cells.Item($newRow, 9).Value  = "Friend"      # RelationValue
$ws.Cells.Item($newRow, 10).Value = "模型高度"     # Desc

# Match the formatting already used on column I ("RelationValue") for the
# preceding rows in this block (text number format).
$ws.Range("I44").NumberFormat = "@"

$ws.Range("J44").Select()
